# Fruta / hortaliza, semanal
# Insert two new weekly price rows into the Repollo sheet, just above the
# existing row 385, pushing the old rows 385-402 down to 387-404.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 385 (old rows 385-402 shift down to 387-404)
$ws.Rows.Item(385).Resize(2).Insert()

# New row 385: Crespo record / Primera
$ws.Cells.Item(385, 1).Value = 7
$ws.Cells.Item(385, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(385, 3).Value = "Ñuble"
$ws.Cells.Item(385, 4).Value = [DateTime]"2023-05-29"
$ws.Cells.Item(385, 5).Value = 16
$ws.Cells.Item(385, 6).Value = 100112006
$ws.Cells.Item(385, 7).Value = "Repollo"
$ws.Cells.Item(385, 8).Value = "Crespo record"
$ws.Cells.Item(385, 9).Value = "Primera"
$ws.Cells.Item(385, 10).Value = 250
$ws.Cells.Item(385, 11).Value = 1200
$ws.Cells.Item(385, 12).Value = 1300
$ws.Cells.Item(385, 13).Value = 1260
$ws.Cells.Item(385, 14).Value = "`$/unidad"
$ws.Cells.Item(385, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(385, 16).Value = 1260
$ws.Cells.Item(385, 17).Value = 1
$ws.Cells.Item(385, 18).Value = "Hortaliza"

# New row 386: Crespo record / Segunda
$ws.Cells.Item(386, 1).Value = 7
$ws.Cells.Item(386, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(386, 3).Value = "Ñuble"
$ws.Cells.Item(386, 4).Value = [DateTime]"2023-05-29"
$ws.Cells.Item(386, 5).Value = 16
$ws.Cells.Item(386, 6).Value = 100112006
$ws.Cells.Item(386, 7).Value = "Repollo"
$ws.Cells.Item(386, 8).Value = "Crespo record"
$ws.Cells.Item(386, 9).Value = "Segunda"
$ws.Cells.Item(386, 10).Value = 120
$ws.Cells.Item(386, 11).Value = 1000
$ws.Cells.Item(386, 12).Value = 1000
$ws.Cells.Item(386, 13).Value = 1000
$ws.Cells.Item(386, 14).Value = "`$/unidad"
$ws.Cells.Item(386, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(386, 16).Value = 1000
$ws.Cells.Item(386, 17).Value = 1
$ws.Cells.Item(386, 18).Value = "Hortaliza"

# Apply the date number format (style) used by the rest of column D to the
# two new date cells, matching the existing D-column formatting.
$ws.Range("D385:D386").NumberFormat = "YYYY-MM-DD HH:MM:SS"
